$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.4578467451490837
$ws.Range("C2").Value = 0.1991884356940687
$ws.Range("D2").Value = 0.03206244218753085
$ws.Range("F2").Value = 0.6078258541310433
$ws.Range("G2").Value = 0.4480220758028679
$ws.Range("H2").Value = 0.6081854275185492
$ws.Range("I2").Value = 0.6185071298194167
$ws.Range("K2").Value = 0.2683249136942436
$ws.Range("L2").Value = 0.3012598701712932
$ws.Range("N2").Value = 1.358636737792086
$ws.Range("O2").Value = 2.066951527494737

$ws.Range("B3").Value = 0.4156696696697963
$ws.Range("C3").Value = 0.2000140003190509
$ws.Range("D3").Value = 0.02915041203873869
$ws.Range("F3").Value = 0.6077144828079994
$ws.Range("G3").Value = 0.4493582754859062
$ws.Range("H3").Value = 0.6119266063073496
$ws.Range("I3").Value = 0.6232048119579048
$ws.Range("K3").Value = 0.2341596531124424
$ws.Range("L3").Value = 0.2898280985136523
$ws.Range("N3").Value = 1.371591582938851
$ws.Range("O3").Value = 2.0773960803173

$ws.Range("B4").Value = 0.389833631767516
$ws.Range("C4").Value = 0.2005562415983206
$ws.Range("D4").Value = 0.0273480889464679
$ws.Range("F4").Value = 0.6079961524622348
$ws.Range("G4").Value = 0.4504627689000529
$ws.Range("H4").Value = 0.6144610090860567
$ws.Range("I4").Value = 0.626390741858593
$ws.Range("K4").Value = 0.2131160319316336
$ws.Range("L4").Value = 0.2829770113289101
$ws.Range("N4").Value = 1.37998325431121
$ws.Range("O4").Value = 2.084900293518373

$ws.Range("B5").Value = 0.3793213003085896
$ws.Range("C5").Value = 0.2007861224685037
$ws.Range("D5").Value = 0.02661006224850837
$ws.Range("F5").Value = 0.6081990384237486
$ws.Range("G5").Value = 0.4509842902937109
$ws.Range("H5").Value = 0.615553539748575
$ws.Range("I5").Value = 0.6277648919774848
$ws.Range("K5").Value = 0.2045245566267795
$ws.Range("L5").Value = 0.2802275198206843
$ws.Range("N5").Value = 1.383512941450142
$ws.Range("O5").Value = 2.088232859220696

$ws.Range("B6").Value = 0.3775767282646427
$ws.Range("C6").Value = 0.2008248331198956
$ws.Range("D6").Value = 0.02648729921359916
$ws.Range("F6").Value = 0.6082380510470671
$ws.Range("G6").Value = 0.4510752028539642
$ws.Range("H6").Value = 0.6157385638682698
$ws.Range("I6").Value = 0.6279976512019481
$ws.Range("K6").Value = 0.2030969973867371
$ws.Range("L6").Value = 0.2797735321612151
$ws.Range("N6").Value = 1.384105689006285
$ws.Range("O6").Value = 2.088802814831126

$ws.Range("B7").Value = 0.3896917927388017
$ws.Range("C7").Value = 0.2005593057274559
$ws.Range("D7").Value = 0.02733815004867779
$ws.Range("F7").Value = 0.6079985318205132
$ws.Range("G7").Value = 0.4504695130949585
$ws.Range("H7").Value = 0.6144755013562815
$ws.Range("I7").Value = 0.6264089669360864
$ws.Range("K7").Value = 0.2130002284355612
$ws.Range("L7").Value = 0.2829397589615326
$ws.Range("N7").Value = 1.380030411414898
$ws.Range("O7").Value = 2.084944125881719

$ws.Range("B8").Value = 0.4432919019999986
$ws.Range("C8").Value = 0.1994657710755732
$ws.Range("D8").Value = 0.03106136855066666
$ws.Range("F8").Value = 0.6077148462408886
$ws.Range("G8").Value = 0.4484238406876599
$ws.Range("H8").Value = 0.6094261747355816
$ws.Range("I8").Value = 0.620064327846741
$ws.Range("K8").Value = 0.256558765160122
$ws.Range("L8").Value = 0.297283381272166
$ws.Range("N8").Value = 1.3630128221007
$ws.Range("O8").Value = 2.070326410916493

$ws.Range("B9").Value = 0.5488561545342634
$ws.Range("C9").Value = 0.1976005997897019
$ws.Range("D9").Value = 0.03824762428452999
$ws.Range("F9").Value = 0.609933350646557
$ws.Range("G9").Value = 0.4466666074002532
$ws.Range("H9").Value = 0.6014044368304212
$ws.Range("I9").Value = 0.6100136693667544
$ws.Range("K9").Value = 0.3414319869407336
$ws.Range("L9").Value = 0.3267411662053092
$ws.Range("N9").Value = 1.333109350813064
$ws.Range("O9").Value = 2.050314443743289

$ws.Range("B10").Value = 0.6266609389056725
$ws.Range("C10").Value = 0.1963988995722659
$ws.Range("D10").Value = 0.0434559831334127
$ws.Range("F10").Value = 0.6132524631356659
$ws.Range("G10").Value = 0.4467511816644816
$ws.Range("H10").Value = 0.5966531766674095
$ws.Range("I10").Value = 0.6040855254925859
$ws.Range("K10").Value = 0.4034334797326835
$ws.Range("L10").Value = 0.3491923933680425
$ws.Range("N10").Value = 1.313249887780145
$ws.Range("O10").Value = 2.040882183414254

$ws.Range("B11").Value = 0.6621039079178672
$ws.Range("C11").Value = 0.195888497725285
$ws.Range("D11").Value = 0.04580964818150335
$ws.Range("F11").Value = 0.6151287576010986
$ws.Range("G11").Value = 0.4470886823423257
$ws.Range("H11").Value = 0.5947389657229252
$ws.Range("I11").Value = 0.6017045051325418
$ws.Range("K11").Value = 0.431558034991383
$ws.Range("L11").Value = 0.3595812859780949
$ws.Range("N11").Value = 1.304672861895462
$ws.Range("O11").Value = 2.037734772874359

$ws.Range("B12").Value = 0.675531628607132
$ws.Range("C12").Value = 0.1957004093787376
$ws.Range("D12").Value = 0.0466986378148988
$ws.Range("F12").Value = 0.6158918893397498
$ws.Range("G12").Value = 0.4472594971434916
$ws.Range("H12").Value = 0.5940495812427145
$ws.Range("I12").Value = 0.6008482478378809
$ws.Range("K12").Value = 0.4421960235016513
$ws.Range("L12").Value = 0.3635404675172111
$ws.Range("N12").Value = 1.301490669980151
$ws.Range("O12").Value = 2.036707243274094

$ws.Range("B13").Value = 0.6726394645887979
$ws.Range("C13").Value = 0.1957406871332807
$ws.Range("D13").Value = 0.04650728062186715
$ws.Range("F13").Value = 0.6157251961159318
$ws.Range("G13").Value = 0.4472207960480006
$ws.Range("H13").Value = 0.5941964751784354
$ws.Range("I13").Value = 0.6010306399993333
$ws.Range("K13").Value = 0.4399054936839093
$ws.Range("L13").Value = 0.362686671279036
$ws.Range("N13").Value = 1.302173087967859
$ws.Range("O13").Value = 2.03692123313067

$ws.Range("B14").Value = 0.663208494355473
$ws.Range("C14").Value = 0.1958729197108724
$ws.Range("D14").Value = 0.04588283207959876
$ws.Range("F14").Value = 0.6151904867356848
$ws.Range("G14").Value = 0.4471018733229215
$ws.Range("H14").Value = 0.5946815387803923
$ws.Range("I14").Value = 0.6016331508760366
$ws.Range("K14").Value = 0.4324334762409876
$ws.Range("L14").Value = 0.3599065072916119
$ws.Range("N14").Value = 1.304409743454961
$ws.Range("O14").Value = 2.0376469445645

$ws.Range("B15").Value = 0.6574325402952752
$ws.Range("C15").Value = 0.195954591130878
$ws.Range("D15").Value = 0.04550003952807913
$ws.Range("F15").Value = 0.6148698123467042
$ws.Range("G15").Value = 0.4470346312793225
$ws.Range("H15").Value = 0.594983273869147
$ws.Range("I15").Value = 0.6020081159403787
$ws.Range("K15").Value = 0.4278550461135069
$ws.Range("L15").Value = 0.3582068456626359
$ws.Range("N15").Value = 1.305788321208286
$ws.Range("O15").Value = 2.038112861258526

$ws.Range("B16").Value = 0.6243455758556991
$ws.Range("C16").Value = 0.1964329830685685
$ws.Range("D16").Value = 0.04330184715217911
$ws.Range("F16").Value = 0.6131372100915087
$ws.Range("G16").Value = 0.4467351433443056
$ws.Range("H16").Value = 0.5967832435946718
$ws.Range("I16").Value = 0.6042474824588204
$ws.Range("K16").Value = 0.4015938053104833
$ws.Range("L16").Value = 0.3485169792154466
$ws.Range("N16").Value = 1.313819615865633
$ws.Range("O16").Value = 2.041110874088574

$ws.Range("B17").Value = 0.6040597978899598
$ws.Range("C17").Value = 0.1967357302948614
$ws.Range("D17").Value = 0.04194929031466188
$ws.Range("F17").Value = 0.6121681125060476
$ws.Range("G17").Value = 0.4466280132822007
$ws.Range("H17").Value = 0.5979507308617471
$ws.Range("I17").Value = 0.6057021132412501
$ws.Range("K17").Value = 0.3854623846487755
$ws.Range("L17").Value = 0.3426174701112643
$ws.Range("N17").Value = 1.318863643287127
$ws.Range("O17").Value = 2.043242840004567

$ws.Range("B18").Value = 0.5923966341687503
$ws.Range("C18").Value = 0.1969132766545059
$ws.Range("D18").Value = 0.04116986551953516
$ws.Range("F18").Value = 0.6116452060412172
$ws.Range("G18").Value = 0.4465945401866662
$ws.Range("H18").Value = 0.5986455059615849
$ws.Range("I18").Value = 0.6065684966748499
$ws.Range("K18").Value = 0.3761765038876774
$ws.Range("L18").Value = 0.3392407759479994
$ws.Range("N18").Value = 1.321807862525606
$ws.Range("O18").Value = 2.044576715416497

$ws.Range("B19").Value = 0.5884485164209821
$ws.Range("C19").Value = 0.1969739779233777
$ws.Range("D19").Value = 0.04090571452190517
$ws.Range("F19").Value = 0.6114740851924481
$ws.Range("G19").Value = 0.4465880404070717
$ws.Range("H19").Value = 0.5988847427575905
$ws.Range("I19").Value = 0.6068669438772858
$ws.Range("K19").Value = 0.3730311929044774
$ws.Range("L19").Value = 0.3381003314325142
$ws.Range("N19").Value = 1.322812114568027
$ws.Range("O19").Value = 2.045046830098741

$ws.Range("B20").Value = 0.6062187753477133
$ws.Range("C20").Value = 0.1967031491585907
$ws.Range("D20").Value = 0.04209342473484412
$ws.Range("F20").Value = 0.6122677052239709
$ws.Range("G20").Value = 0.4466365044478735
$ws.Range("H20").Value = 0.5978240421749774
$ws.Range("I20").Value = 0.605544189703096
$ws.Range("K20").Value = 0.3871803841944939
$ws.Range("L20").Value = 0.3432437712986456
$ws.Range("N20").Value = 1.318322244683159
$ws.Range("O20").Value = 2.043004750448659

$ws.Range("B21").Value = 0.6659784362949495
$ws.Range("C21").Value = 0.1958339391395718
$ws.Range("D21").Value = 0.04606631030557651
$ws.Range("F21").Value = 0.6153461163213336
$ws.Range("G21").Value = 0.4471356364322361
$ws.Range("H21").Value = 0.5945381012159459
$ws.Range("I21").Value = 0.6014549472695876
$ws.Range("K21").Value = 0.4346285242257864
$ws.Range("L21").Value = 0.3607224278142667
$ws.Range("N21").Value = 1.303750999262309
$ws.Range("O21").Value = 2.037429326469635

$ws.Range("B22").Value = 0.7050707984833764
$ws.Range("C22").Value = 0.1952960989958576
$ws.Range("D22").Value = 0.04864943980690839
$ws.Range("F22").Value = 0.6176647107858528
$ws.Range("G22").Value = 0.4477125512461271
$ws.Range("H22").Value = 0.5925973606305774
$ws.Range("I22").Value = 0.5990469123200484
$ws.Range("K22").Value = 0.4655673957764179
$ws.Range("L22").Value = 0.3722921732435083
$ws.Range("N22").Value = 1.294611044400476
$ws.Range("O22").Value = 2.034743242506721

$ws.Range("B23").Value = 0.6842034539228052
$ws.Range("C23").Value = 0.1955803954366608
$ws.Range("D23").Value = 0.04727201349494692
$ws.Range("F23").Value = 0.6163991959490076
$ws.Range("G23").Value = 0.4473816982383738
$ws.Range("H23").Value = 0.5936142655181129
$ws.Range("I23").Value = 0.6003079278447103
$ws.Range("K23").Value = 0.4490614644408879
$ws.Range("L23").Value = 0.366103828590937
$ws.Range("N23").Value = 1.299454145657798
$ws.Range("O23").Value = 2.036089249165883

$ws.Range("B24").Value = 0.6052427034435084
$ws.Range("C24").Value = 0.1967178682087685
$ws.Range("D24").Value = 0.04202826723494013
$ws.Range("F24").Value = 0.6122225726948329
$ws.Range("G24").Value = 0.4466325780098117
$ws.Range("H24").Value = 0.5978812447039701
$ws.Range("I24").Value = 0.6056154931844482
$ws.Range("K24").Value = 0.3864037130339852
$ws.Range("L24").Value = 0.3429605737614168
$ws.Range("N24").Value = 1.318566872873966
$ws.Range("O24").Value = 2.043112053759529

$ws.Range("B25").Value = 0.520252700251234
$ws.Range("C25").Value = 0.1980754484435359
$ws.Range("D25").Value = 0.03631599290952892
$ws.Range("F25").Value = 0.6090363708276598
$ws.Range("G25").Value = 0.4469005058674398
$ws.Range("H25").Value = 0.6033736559278182
$ws.Range("I25").Value = 0.6124768178799727
$ws.Range("K25").Value = 0.3185322123113963
$ws.Range("L25").Value = 0.3186298929525719
$ws.Range("N25").Value = 1.340828137614245
$ws.Range("O25").Value = 2.05480230170798

